$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

$ws.Range("D2").Value = "26.369.43"

$ws.Range("D3").Value = "1.679.10"
$ws.Range("E3").Value = "  +2.05%  "

Set-TextValue $ws.Range("D4") "0.9987"
$ws.Range("E4").Value = "  -0.41%  "

Set-TextValue $ws.Range("D5") "217.70"
$ws.Range("E5").Value = "  +5.35%  "

Set-TextValue $ws.Range("D6") "0.5309"
$ws.Range("E6").Value = "  +2.28%  "

Set-TextValue $ws.Range("D7") "0.9993"
$ws.Range("E7").Value = "  -0.42%  "

Set-TextValue $ws.Range("D8") "0.2659"
$ws.Range("E8").Value = "  +3.30%  "

Set-TextValue $ws.Range("D9") "0.06454"
$ws.Range("E9").Value = "  +3.29%  "

Set-TextValue $ws.Range("D10") "21.22"
$ws.Range("E10").Value = "  +2.36%  "

Set-TextValue $ws.Range("D11") "0.07795"
$ws.Range("E11").Value = "  +3.04%  "

$ws.Range("D12").Value = "1.678.11"
$ws.Range("E12").Value = "  +2.24%  "

Set-TextValue $ws.Range("D13") "4.505"
$ws.Range("E13").Value = "  +2.81%  "

$ws.Range("D14").Value = "1.901.17"
$ws.Range("E14").Value = "  +1.78%  "

Set-TextValue $ws.Range("D15") "0.5605"
$ws.Range("E15").Value = "  +3.95%  "

$ws.Range("D16").Value = "0.0₅8442"
$ws.Range("E16").Value = "  +6.51%  "

Set-TextValue $ws.Range("D17") "66.09"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "26.379.21"
$ws.Range("E18").Value = "  +1.49%  "

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "4.833"
$ws.Range("E19").Value = "  +3.49%  "

$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D20") "0.9998"
$ws.Range("E20").Value = "  -0.53%  "

Set-TextValue $ws.Range("D21") "195.50"
$ws.Range("E21").Value = "  +4.32%  "

Set-TextValue $ws.Range("D22") "10.39"
$ws.Range("E22").Value = "  +3.91%  "

Set-TextValue $ws.Range("D23") "6.402"
$ws.Range("E23").Value = "  +4.50%  "

Set-TextValue $ws.Range("D24") "1.0000"
$ws.Range("E24").Value = "  -0.38%  "

Set-TextValue $ws.Range("D25") "142.97"
$ws.Range("E25").Value = "  -3.57%  "

$ws.Range("E26").Value = "  +4.13%  "

Set-TextValue $ws.Range("D27") "7.486"
$ws.Range("E27").Value = "  +1.85%  "

Set-TextValue $ws.Range("D28") "16.31"
$ws.Range("E28").Value = "  +4.36%  "

Set-TextValue $ws.Range("D29") "1.435"
$ws.Range("E29").Value = "  +3.50%  "

Set-TextValue $ws.Range("D30") "0.06203"
$ws.Range("E30").Value = "  +3.18%  "

Set-TextValue $ws.Range("D31") "1.276"
$ws.Range("E31").Value = "  +2.84%  "

Set-TextValue $ws.Range("D32") "3.545"
$ws.Range("E32").Value = "  +3.01%  "

Set-TextValue $ws.Range("D33") "3.461"
$ws.Range("E33").Value = "  +1.88%  "

Set-TextValue $ws.Range("D34") "1.708"
$ws.Range("E34").Value = "  +4.97%  "

Set-TextValue $ws.Range("D35") "1.014"
$ws.Range("E35").Value = "  +3.84%  "

Set-TextValue $ws.Range("D36") "2.782"
$ws.Range("E36").Value = "  +2.14%  "

Set-TextValue $ws.Range("D37") "2.402"
$ws.Range("E37").Value = "  +0.79%  "

Set-TextValue $ws.Range("D38") "0.5747"
$ws.Range("E38").Value = "  -1.80%  "

Set-TextValue $ws.Range("D39") "0.01633"
$ws.Range("E39").Value = "  +2.99%  "

Set-TextValue $ws.Range("D40") "5.936"
$ws.Range("E40").Value = "  -0.53%  "

Set-TextValue $ws.Range("D41") "0.8662"
$ws.Range("E41").Value = "  +2.67%  "

$ws.Range("D42").Value = "1.054.54"
$ws.Range("E42").Value = "  -2.87%  "

Set-TextValue $ws.Range("D43") "0.9992"
$ws.Range("E43").Value = "  -0.38%  "

Set-TextValue $ws.Range("D44") "100.05"
$ws.Range("E44").Value = "  -0.28%  "

$ws.Range("D45").Value = "1.824.83"
$ws.Range("E45").Value = "  +1.59%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈110"
$ws.Range("E46").Value = "  +5.13%  "

$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D47") "57.03"
$ws.Range("E47").Value = "  +4.06%  "

Set-TextValue $ws.Range("D48") "8.160"
$ws.Range("E48").Value = "  +2.41%  "

$ws.Range("E49").Value = "  -0.31%  "

Set-TextValue $ws.Range("D50") "0.05193"
$ws.Range("E50").Value = "  -0.54%  "

Set-TextValue $ws.Range("D51") "6.065"
$ws.Range("E51").Value = "  +3.53%  "
